$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")
Write-Host "hello"
